$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 9.295069333333332
$ws.Range("H2").Value = 27.885208
$ws.Range("I2").Value = 0.2851098797714356
$ws.Range("J2").Value = 0.2851098797714357
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.044118333333333
$ws.Range("N2").Value = 6.132354999999999
$ws.Range("O2").Value = 0.1776005292722278
$ws.Range("P2").Value = 0.1776005292722278
$ws.Range("Q2").Value = 19.00022163387111
$ws.Range("R2").Value = 171.0019947048399
$ws.Range("S2").Value = 0.0506356655481482
$ws.Range("T2").Value = 0.05063566554814821

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 9.295069333333332
$ws.Range("H3").Value = 27.885208
$ws.Range("I3").Value = 0.2851098797714356
$ws.Range("J3").Value = 0.2851098797714357
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.059280333333334
$ws.Range("N3").Value = 21.177841
$ws.Range("O3").Value = 0.6133362746356149
$ws.Range("P3").Value = 0.6133362746356149
$ws.Range("Q3").Value = 65.61650014176978
$ws.Range("R3").Value = 590.548501275928
$ws.Range("S3").Value = 0.1748682315208204
$ws.Range("T3").Value = 0.1748682315208204

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 9.295069333333332
$ws.Range("H4").Value = 27.885208
$ws.Range("I4").Value = 0.2851098797714356
$ws.Range("J4").Value = 0.2851098797714357
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 2.406242333333333
$ws.Range("N4").Value = 7.218726999999999
$ws.Range("O4").Value = 0.2090631960921573
$ws.Range("P4").Value = 0.2090631960921573
$ws.Range("Q4").Value = 22.36618932113511
$ws.Range("R4").Value = 201.295703890216
$ws.Range("S4").Value = 0.05960598270246703
$ws.Range("T4").Value = 0.05960598270246705

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 8.588082333333332
$ws.Range("H5").Value = 25.764247
$ws.Range("I5").Value = 0.2634242988100204
$ws.Range("J5").Value = 0.2634242988100204
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.044118333333333
$ws.Range("N5").Value = 6.132354999999999
$ws.Range("O5").Value = 0.1776005292722278
$ws.Range("P5").Value = 0.1776005292722278
$ws.Range("Q5").Value = 17.55505654574277
$ws.Range("R5").Value = 157.9955089116849
$ws.Range("S5").Value = 0.04678429489182511
$ws.Range("T5").Value = 0.04678429489182511

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 8.588082333333332
$ws.Range("H6").Value = 25.764247
$ws.Range("I6").Value = 0.2634242988100204
$ws.Range("J6").Value = 0.2634242988100204
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 7.059280333333334
$ws.Range("N6").Value = 21.177841
$ws.Range("O6").Value = 0.6133362746356149
$ws.Range("P6").Value = 0.6133362746356149
$ws.Range("Q6").Value = 60.62568071674744
$ws.Range("R6").Value = 545.6311264507269
$ws.Range("S6").Value = 0.161567678080637
$ws.Range("T6").Value = 0.161567678080637

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 8.588082333333332
$ws.Range("H7").Value = 25.764247
$ws.Range("I7").Value = 0.2634242988100204
$ws.Range("J7").Value = 0.2634242988100204
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.406242333333333
$ws.Range("N7").Value = 7.218726999999999
$ws.Range("O7").Value = 0.2090631960921573
$ws.Range("P7").Value = 0.2090631960921573
$ws.Range("Q7").Value = 20.66500727261877
$ws.Range("R7").Value = 185.985065453569
$ws.Range("S7").Value = 0.05507232583755833
$ws.Range("T7").Value = 0.05507232583755834

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.71855733333333
$ws.Range("H8").Value = 44.155672
$ws.Range("I8").Value = 0.4514658214185439
$ws.Range("J8").Value = 0.4514658214185439
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.044118333333333
$ws.Range("N8").Value = 6.132354999999999
$ws.Range("O8").Value = 0.1776005292722278
$ws.Range("P8").Value = 0.1776005292722278
$ws.Range("Q8").Value = 30.08647288528444
$ws.Range("R8").Value = 270.77825596756
$ws.Range("S8").Value = 0.08018056883225447
$ws.Range("T8").Value = 0.08018056883225445

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.71855733333333
$ws.Range("H9").Value = 44.155672
$ws.Range("I9").Value = 0.4514658214185439
$ws.Range("J9").Value = 0.4514658214185439
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.059280333333334
$ws.Range("N9").Value = 21.177841
$ws.Range("O9").Value = 0.6133362746356149
$ws.Range("P9").Value = 0.6133362746356149
$ws.Range("Q9").Value = 103.9024223182391
$ws.Range("R9").Value = 935.1218008641521
$ws.Range("S9").Value = 0.2769003650341575
$ws.Range("T9").Value = 0.2769003650341575

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.71855733333333
$ws.Range("H10").Value = 44.155672
$ws.Range("I10").Value = 0.4514658214185439
$ws.Range("J10").Value = 0.4514658214185439
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 2.406242333333333
$ws.Range("N10").Value = 7.218726999999999
$ws.Range("O10").Value = 0.2090631960921573
$ws.Range("P10").Value = 0.2090631960921573
$ws.Range("Q10").Value = 35.41641574106045
$ws.Range("R10").Value = 318.747741669544
$ws.Range("S10").Value = 0.0943848875521319
$ws.Range("T10").Value = 0.0943848875521319
